# Update the marksheet with corrected correct/total marks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row, Right column: 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row, Right column: 54 -> 90
$ws.Range("B12").Value = 90

# "Total" row, Max column text: "48/84" -> "90/140"
$ws.Range("E12").Value = "90/140"
